# Weekly update: insert a new price observation at row 83 for
# "Hortaliza, Terminal La Palmera de La Serena - Ajo" and push the
# existing rows (old 83..207) down by one (new 84..208).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 83; Excel shifts rows
# 83..207 down to 84..208 automatically (and copies formatting from the
# row above, matching native Excel "Insert" behaviour).
$ws.Rows(83).Insert()

# Populate the newly inserted row 83 with the latest observation.
$ws.Range("A83").Value = 8
$ws.Range("B83").Value = "Terminal La Palmera de La Serena"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = 44580
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = 100112003
$ws.Range("G83").Value = "Ajo"
$ws.Range("H83").Value = "Chino"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 19000
$ws.Range("M83").Value = 18500
$ws.Range("N83").Value = "`$/caja 10 kilos"
$ws.Range("O83").Value = "China"
$ws.Range("P83").Value = 1850
$ws.Range("Q83").Value = 10
$ws.Range("R83").Value = "Hortaliza"
